# Update cryptos list (prices, 1h volume %, and a few coin re-rankings)
# matching the "Updated cryptos list ... with GitHub Actions" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.164.30"
$ws.Range("E2").Value = "  +0.48%  "
$ws.Range("D3").Value = "1.856.14"
$ws.Range("E3").Value = "  +1.75%  "
$ws.Range("E4").Value = "  +0.47%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.29"
$ws.Range("E5").Value = "  +3.21%  "
$ws.Range("E6").Value = "  +0.52%  "
$ws.Range("E7").Value = "  +0.41%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "41.87"
$ws.Range("E8").Value = "  +4.74%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.327"
$ws.Range("E9").Value = "  +1.67%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0693"
$ws.Range("E10").Value = "  +1.40%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0990"
$ws.Range("E11").Value = "  -0.27%  "
$ws.Range("D12").Value = "2.123.95"
$ws.Range("E12").Value = "  +1.67%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.41"
$ws.Range("E13").Value = "  +0.66%  "
$ws.Range("D14").Value = "1.845.74"
$ws.Range("E14").Value = "  +1.16%  "
$ws.Range("E15").Value = "  +1.13%  "
$ws.Range("E16").Value = "  +1.25%  "
$ws.Range("D17").Value = "35.142.80"
$ws.Range("E17").Value = "  +0.36%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.92"
$ws.Range("E18").Value = "  +0.30%  "
$ws.Range("E19").Value = "  +0.47%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "240.77"
$ws.Range("E20").Value = "  +0.11%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.20"
$ws.Range("E21").Value = "  +0.76%  "
$ws.Range("E22").Value = "  +1.04%  "
$ws.Range("E23").Value = "  +0.44%  "
$ws.Range("E24").Value = "  +0.19%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "167.55"
$ws.Range("E25").Value = "  -3.52%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.86"
$ws.Range("E26").Value = "  +22.33%  "
$ws.Range("E27").Value = "  +1.71%  "
$ws.Range("E28").Value = "  +1.33%  "
$ws.Range("E29").Value = "  -1.05%  "
$ws.Range("E30").Value = "  +0.44%  "
$ws.Range("E31").Value = "  +0.63%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.98"
$ws.Range("E32").Value = "  -0.18%  "
$ws.Range("B33").Value = "WEMIXToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.80"
$ws.Range("E33").Value = "  +28.12%  "
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.98"
$ws.Range("E34").Value = "  +0.52%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.847"
$ws.Range("E35").Value = "  +21.35%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.01"
$ws.Range("E36").Value = "  +10.42%  "
$ws.Range("E37").Value = "  +4.89%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.06"
$ws.Range("E38").Value = "  +6.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "89.99"
$ws.Range("E39").Value = "  -3.22%  "
$ws.Range("E40").Value = "  +3.34%  "
$ws.Range("D41").Value = "1.340.55"
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("E42").Value = "  +0.97%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.31"
$ws.Range("E43").Value = "  +2.23%  "
$ws.Range("B44").Value = "HuobiToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.42"
$ws.Range("E44").Value = "  -0.11%  "
$ws.Range("B45").Value = "Kaspa"
$ws.Range("C45").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0557"
$ws.Range("E45").Value = "  +6.61%  "
$ws.Range("B46").Value = "MXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.74"
$ws.Range("E46").Value = "  -0.47%  "
$ws.Range("B47").Value = "Gas"
$ws.Range("C47").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.17"
$ws.Range("E47").Value = "  +43.76%  "
$ws.Range("E48").Value = "  +3.03%  "
$ws.Range("D49").Value = "2.034.03"
$ws.Range("E49").Value = "  +1.25%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0678"
$ws.Range("E50").Value = "  +1.08%  "
$ws.Range("E51").Value = "  +0.42%  "
